# The two header cells ("Month" and "Savings") keep the same Bold/Italic/
# StrikeThrough formatting (bold on, italic off, strike off) - this commit
# only changed how the boolean w:val was serialized (POI 4.1.0 "true"/
# "false" -> POI 5.2.3 "on"/"off"). Re-apply the same formatting values so
# the run properties get rewritten with this runtime's canonical boolean
# form.
$d = $word.ActiveDocument

$targets = @("Month", "Savings")
foreach ($t in $targets) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Font.Italic = 0
        $rng.Font.StrikeThrough = 0
    }
}
